$d = $word.ActiveDocument

# 1) "-" + bookmark + "Amministratore:" -> "-Gestione Amministratore:"
#    (this also removes the pre-existing "_GoBack" bookmark that sat between
#    the "-" run and the "Amministratore:" run)
$d.Content.Find.Execute("-Amministratore:", $true, $false, $false, $false, $false, $true, 1, $false, "-Gestione Amministratore:", 2)

# 2) "- Ricerca:" -> "- Gestione Ricerca:" (the section heading, not the
#    pre-existing "- Gestione Ricerca;" bullet earlier in the document)
$d.Content.Find.Execute("- Ricerca:", $true, $false, $false, $false, $false, $true, 1, $false, "- Gestione Ricerca:", 2)

# 3) Re-add the "_GoBack" bookmark right after the "8 Test Cases" heading run.
#    A bookmark collapsed exactly at a paragraph's end boundary lands in the
#    wrong spot, so nudge it off that boundary with a throwaway character,
#    add the bookmark, then remove the throwaway character again.
$range = $d.Content
$range.Find.Execute("8 Test Cases", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $range.End

$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bookmarkRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$d.Range($endPos, $endPos + 1).Delete()
